$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "last changed by / status" columns (C1 "aanpassen" header / D1 "aangepast"
# header) are removed - the sheet goes back to just two columns (use case + person).
$ws.Range("C1").ClearContents()
$ws.Range("D1").ClearContents()

# "zoek functie" (row 12) now also has an owner, just like the other use cases: Eric.
$ws.Range("B12").Value = "Eric"

# Column A is widened to fit the longer text, and the last touched cell becomes B14.
$ws.Columns("A").ColumnWidth = 41.17
$ws.Range("B14").Select()
